# Generate Report for Handoff
# Refresh the "Latest Handoff" timestamp for the 159503f2-24c2-4ef9-ae22-14852ef12a3c
# row (row 5) across the Overview / zh-cn / de-de status sheets, exactly like a fresh
# handoff-report run would: the row was just (re-)handed off, so its recorded
# handoff date/time moves forward a few seconds/minutes.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-22 04:38:10"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-22 04:38:06"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-22 04:38:10"
